# Updates the cryptos price/volume table to the latest scraped values
# (GitHub Actions refresh), including two rows whose coin/link/price/volume
# swapped ranking position with their neighbour.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values (column D) look like plain decimals (e.g. "306.06").
# If assigned directly, Excel would auto-convert them into numbers instead
# of keeping them as text. Temporarily force a text number format, assign
# the value, then restore the cell's original style so no visible
# formatting/style changes are left behind.
function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "42.915.86"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.302.35"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "306.06"
$ws.Range("E5").Value = "  +1.75%  "
Set-TextValue "D6" "97.43"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.35%  "
Set-TextValue "D10" "35.75"
$ws.Range("E10").Value = "  -0.27%  "
Set-TextValue "D11" "0.0791"
$ws.Range("E11").Value = "  +0.07%  "
Set-TextValue "D12" "18.19"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("E13").Value = "  +0.88%  "
Set-TextValue "D14" "6.78"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "2.660.92"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "2.302.38"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "42.846.13"
$ws.Range("E18").Value = "  -0.46%  "
Set-TextValue "D19" "12.69"
$ws.Range("E19").Value = "  -4.80%  "
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  -1.09%  "
Set-TextValue "D22" "67.77"
$ws.Range("E22").Value = "  -1.25%  "
Set-TextValue "D23" "236.76"
$ws.Range("E23").Value = "  -0.55%  "
Set-TextValue "D24" "2.15"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -0.01%  "
Set-TextValue "D28" "25.47"
$ws.Range("E28").Value = "  +2.65%  "
Set-TextValue "D29" "166.98"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("E30").Value = "  +1.51%  "
Set-TextValue "D31" "9.06"
$ws.Range("E31").Value = "  -1.10%  "
Set-TextValue "D32" "33.10"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("E33").Value = "  +0.09%  "
Set-TextValue "D34" "4.82"
$ws.Range("E34").Value = "  -0.29%  "
Set-TextValue "D36" "17.25"
$ws.Range("E36").Value = "  -4.67%  "
Set-TextValue "D37" "2.39"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("E41").Value = "  -1.26%  "
Set-TextValue "D42" "2.73"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").Value = "2.007.27"
$ws.Range("E43").Value = "  +0.12%  "
Set-TextValue "D44" "0.0281"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "17.99"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D47" "2.12"
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D49" "2.88"
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D50" "54.15"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "2.528.58"
$ws.Range("E51").Value = "  +0.01%  "
